$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-09 Saturday" "2024-11-10 Sunday"

Replace-Text "63×69=4347" "23×46=1058"
Replace-Text "98×32=3136" "61×44=2684"
Replace-Text "24×94=2256" "28×31=868"
Replace-Text "59×97=5723" "98×28=2744"
Replace-Text "63×84=5292" "35×76=2660"

Replace-Text "34×57=1938" "12×56=672"
Replace-Text "53×69=3657" "26×87=2262"
Replace-Text "33×78=2574" "68×11=748"
Replace-Text "31×95=2945" "76×23=1748"
Replace-Text "34×18=612" "76×46=3496"

Replace-Text "53×43=2279" "73×91=6643"
Replace-Text "82×65=5330" "36×13=468"
Replace-Text "49×78=3822" "38×54=2052"
Replace-Text "93×58=5394" "98×14=1372"
Replace-Text "87×93=8091" "12×42=504"

Replace-Text "49×26=1274" "67×55=3685"
Replace-Text "48×32=1536" "26×76=1976"
Replace-Text "44×68=2992" "55×13=715"
Replace-Text "73×95=6935" "85×27=2295"
Replace-Text "35×28=980" "76×36=2736"

Replace-Text "95×92=8740" "47×54=2538"
Replace-Text "30×43=1290" "56×77=4312"
Replace-Text "41×16=656" "57×38=2166"
Replace-Text "71×75=5325" "82×29=2378"
Replace-Text "17×11=187" "28×76=2128"
